# "Accept deal alert added"
# Roll the dated test accounts (031816 -> 032816) forward on the four
# Nymgo sign-up sheets, and move the active tab/selection to
# NymgoInterReseller (matching a manual walk-through of the new accounts).

$wb = $excel.ActiveWorkbook

# NymgoEuroNormalUser: edeal031816 -> edeal032816
$ws13 = $wb.Worksheets.Item("NymgoEuroNormalUser")
$ws13.Range("B1").Value = "edeal032816"
$ws13.Range("B3").Value = "edeal032816"
$ws13.Range("B4").Value = "edeal032816@mail.ru"
$ws13.Range("J26").Select() | Out-Null

# NymgoEuroReseller: edealReseller031816 -> edealReseller032816
$ws14 = $wb.Worksheets.Item("NymgoEuroReseller")
$ws14.Range("B1").Value = "edealReseller032816"
$ws14.Range("B3").Value = "edealReseller032816"
$ws14.Range("B4").Value = "edealReseller032816@mail.ru"
$ws14.Range("D24").Select() | Out-Null

# NymgoInterNormalUser: ideal031816 -> ideal032816
$ws15 = $wb.Worksheets.Item("NymgoInterNormalUser")
$ws15.Range("B1").Value = "ideal032816"
$ws15.Range("B3").Value = "ideal032816"
$ws15.Range("B4").Value = "ideal032816@mail.ru"

# NymgoInterReseller: idealReseller031816 -> idealReseller032816
$ws16 = $wb.Worksheets.Item("NymgoInterReseller")
$ws16.Range("B1").Value = "idealReseller032816"
$ws16.Range("B3").Value = "idealReseller032816"
$ws16.Range("B4").Value = "idealReseller032816@mail.ru"

# Leave the workbook focused on NymgoInterReseller, as in the authored edit.
$ws16.Activate() | Out-Null
$ws16.Range("M18").Select() | Out-Null
